$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 536, pushing existing data (rows 536 onward) down by two rows.
$ws.Rows("536:537").Insert()

# --- New row 536 ---
$ws.Range("A536").Value = 7
$ws.Range("B536").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C536").Value = "Ñuble"
$ws.Range("D536").Value = 44785
$ws.Range("E536").Value = 16
$ws.Range("F536").Value = "Fruta"
$ws.Range("G536").Value = 100108
$ws.Range("H536").Value = "Tropicales y subtropicales"
$ws.Range("I536").Value = 100108006
$ws.Range("J536").Value = "Plátano"
$ws.Range("K536").Value = "Sin especificar"
$ws.Range("L536").Value = "Pintón"
$ws.Range("M536").Value = 80
$ws.Range("N536").Value = 21000
$ws.Range("O536").Value = 21000
$ws.Range("P536").Value = 21000
$ws.Range("Q536").Value = "$/caja 20 kilos"
$ws.Range("R536").Value = "Ecuador"
$ws.Range("S536").Value = 1050
$ws.Range("T536").Value = 20

# --- New row 537 ---
$ws.Range("A537").Value = 7
$ws.Range("B537").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C537").Value = "Ñuble"
$ws.Range("D537").Value = 44785
$ws.Range("E537").Value = 16
$ws.Range("F537").Value = "Fruta"
$ws.Range("G537").Value = 100108
$ws.Range("H537").Value = "Tropicales y subtropicales"
$ws.Range("I537").Value = 100108006
$ws.Range("J537").Value = "Plátano"
$ws.Range("K537").Value = "Sin especificar"
$ws.Range("L537").Value = "Primera Pintón"
$ws.Range("M537").Value = 160
$ws.Range("N537").Value = 22000
$ws.Range("O537").Value = 23000
$ws.Range("P537").Value = 22500
$ws.Range("Q537").Value = "$/caja 20 kilos"
$ws.Range("R537").Value = "Ecuador"
$ws.Range("S537").Value = 1125
$ws.Range("T537").Value = 20
